$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

# Append 17 new OPT-sensor recording rows (995-1011) below the last
# existing row (994). Each new row is first cloned from row 994 via
# copy/paste-special so the numeric-looking text values in columns A
# ("025") and J ("1") stay text-typed (matching the existing data) rather
# than being auto-coerced to numbers by a plain .Value assignment.
for ($r = 995; $r -le 1011; $r++) {
    $ws.Range("A994:J994").Copy()
    $ws.Range("A" + $r + ":J" + $r).PasteSpecial(-4104)
}
$excel.CutCopyMode = 0

# Now fill in the per-row folder (C) / recording (D) values, and the
# protocol (G) for row 996 which differs from the rest.
$ws.Range("C995").Value = "BL-003_025_230619"
$ws.Range("D995").Value = "230619_235230"
$ws.Range("C996").Value = "BL-003_025_230619"
$ws.Range("D996").Value = "230619_235558"
$ws.Range("G996").Value = "Exp#1_60s"
$ws.Range("C997").Value = "BL-003_025_230620"
$ws.Range("D997").Value = "230620_002256"
$ws.Range("C998").Value = "BL-003_025_230620"
$ws.Range("D998").Value = "230620_002719"
$ws.Range("C999").Value = "BL-003_025_230620"
$ws.Range("D999").Value = "230620_003239"
$ws.Range("C1000").Value = "BL-003_025_230620"
$ws.Range("D1000").Value = "230620_003430"
$ws.Range("C1001").Value = "BL-003_025_230620"
$ws.Range("D1001").Value = "230620_003443"
$ws.Range("C1002").Value = "BL-003_025_230620"
$ws.Range("D1002").Value = "230620_003713"
$ws.Range("C1003").Value = "BL-003_025_230620"
$ws.Range("D1003").Value = "230620_003756"
$ws.Range("C1004").Value = "BL-003_025_230620"
$ws.Range("D1004").Value = "230620_004343"
$ws.Range("C1005").Value = "BL-003_025_230620"
$ws.Range("D1005").Value = "230620_005035"
$ws.Range("C1006").Value = "BL-003_025_230620"
$ws.Range("D1006").Value = "230620_005851"
$ws.Range("C1007").Value = "BL-003_025_230620"
$ws.Range("D1007").Value = "230620_010150"
$ws.Range("C1008").Value = "BL-003_025_230620"
$ws.Range("D1008").Value = "230620_011017"
$ws.Range("C1009").Value = "BL-003_025_230620"
$ws.Range("D1009").Value = "230620_012151"
$ws.Range("C1010").Value = "BL-003_025_230620"
$ws.Range("D1010").Value = "230620_012558"
$ws.Range("C1011").Value = "BL-003_025_230620"
$ws.Range("D1011").Value = "230620_012711"

# Match the workbook-level active-tab change: focus moves to DB.
$ws.Activate()
